# Auto-generated PowerShell Excel COM-interop script
# Applies odds/value updates to rows 2,4,5,6,7,8,10,11,13,14,15

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.76
$ws.Range("G2").Value = 1.9
$ws.Range("H2").Value = 4.6
$ws.Range("I2").Value = 5.8
$ws.Range("J2").Value = 3.65
$ws.Range("K2").Value = 4.3
$ws.Range("L2").Value = 1.31
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 3.8
$ws.Range("O2").Value = 1.31
$ws.Range("P2").Value = 1.97
$ws.Range("Q2").Value = 1.83
$ws.Range("R2").Value = 1.37
$ws.Range("S2").Value = 3.1
$ws.Range("T2").Value = 1.78
$ws.Range("U2").Value = 2.04
$ws.Range("V2").Value = 1.2
$ws.Range("W2").Value = 2.1
$ws.Range("X2").Value = 16.5
$ws.Range("Y2").Value = 21
$ws.Range("Z2").Value = 46
$ws.Range("AA2").Value = 140
$ws.Range("AB2").Value = 9.800000000000001
$ws.Range("AC2").Value = 10.5
$ws.Range("AD2").Value = 24
$ws.Range("AE2").Value = 80
$ws.Range("AF2").Value = 13
$ws.Range("AG2").Value = 12.5
$ws.Range("AH2").Value = 23
$ws.Range("AI2").Value = 80
$ws.Range("AJ2").Value = 24
$ws.Range("AK2").Value = 23
$ws.Range("AL2").Value = 42
$ws.Range("AM2").Value = 130
$ws.Range("AN2").Value = 14
$ws.Range("AO2").Value = 85
# Row 4
$ws.Range("F4").Value = 1.71
$ws.Range("G4").Value = 1.86
$ws.Range("I4").Value = 5.4
$ws.Range("N4").Value = 4.2
$ws.Range("P4").Value = 2.1
$ws.Range("Q4").Value = 1.63
$ws.Range("R4").Value = 1.44
$ws.Range("S4").Value = 2.86
$ws.Range("T4").Value = 1.72
$ws.Range("V4").Value = 1.22
# Row 5
$ws.Range("F5").Value = 2.06
$ws.Range("G5").Value = 2.7
$ws.Range("H5").Value = 2.32
$ws.Range("I5").Value = 3.45
$ws.Range("J5").Value = 3.2
$ws.Range("V5").Value = 1.01
# Row 6
$ws.Range("F6").Value = 2.02
$ws.Range("G6").Value = 2.36
$ws.Range("I6").Value = 4.4
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 3.9
$ws.Range("L6").Value = 1.42
$ws.Range("M6").Value = 1.08
$ws.Range("N6").Value = 2.88
$ws.Range("O6").Value = 1.34
$ws.Range("P6").Value = 1.73
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 1.28
$ws.Range("S6").Value = 3.35
$ws.Range("T6").Value = 1.8
$ws.Range("U6").Value = 1.94
$ws.Range("V6").Value = 1.3
$ws.Range("W6").Value = 1.73
$ws.Range("Y6").Value = 14
$ws.Range("Z6").Value = 30
$ws.Range("AA6").Value = 90
$ws.Range("AB6").Value = 9.4
$ws.Range("AD6").Value = 17.5
$ws.Range("AE6").Value = 55
$ws.Range("AF6").Value = 14.5
$ws.Range("AG6").Value = 12
$ws.Range("AH6").Value = 21
$ws.Range("AI6").Value = 65
$ws.Range("AJ6").Value = 980
$ws.Range("AK6").Value = 27
$ws.Range("AL6").Value = 46
$ws.Range("AN6").Value = 24
$ws.Range("AO6").Value = 70
# Row 7
$ws.Range("F7").Value = 2.64
$ws.Range("H7").Value = 2.88
$ws.Range("K7").Value = 3.35
$ws.Range("L7").Value = 1.48
$ws.Range("N7").Value = 3.1
$ws.Range("O7").Value = 1.41
$ws.Range("P7").Value = 1.71
$ws.Range("Q7").Value = 2.2
$ws.Range("T7").Value = 1.87
$ws.Range("X7").Value = 13.5
$ws.Range("AI7").Value = 55
$ws.Range("AN7").Value = 38
# Row 8
$ws.Range("F8").Value = 3.65
$ws.Range("I8").Value = 2.6
$ws.Range("J8").Value = 2.8
$ws.Range("N8").Value = 2.28
$ws.Range("O8").Value = 1.68
$ws.Range("P8").Value = 1.42
$ws.Range("V8").Value = 1.62
$ws.Range("AI8").Value = 100
# Row 10
$ws.Range("H10").Value = 3.85
$ws.Range("L10").Value = 1.56
$ws.Range("O10").Value = 1.52
# Row 11
$ws.Range("I11").Value = 4.5
$ws.Range("K11").Value = 3.6
$ws.Range("N11").Value = 3.05
$ws.Range("V11").Value = 1.28
# Row 13
$ws.Range("G13").Value = 1.69
$ws.Range("V13").Value = 1.11
# Row 14
$ws.Range("F14").Value = 2.52
$ws.Range("G14").Value = 2.78
$ws.Range("J14").Value = 2.9
$ws.Range("L14").Value = 1.61
$ws.Range("N14").Value = 2.46
# Row 15
$ws.Range("F15").Value = 1.7
$ws.Range("H15").Value = 6.4
$ws.Range("J15").Value = 3.45
$ws.Range("K15").Value = 3.75
$ws.Range("L15").Value = 1.51
$ws.Range("N15").Value = 2.84
$ws.Range("O15").Value = 1.48
